# Adding math-prefix to fillers file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Rename header cells to have a "math_" prefix
$ws.Range("C1").Value = "math_addition"
$ws.Range("D1").Value = "math_subtraction"
$ws.Range("E1").Value = "math_multiplication"
$ws.Range("F1").Value = "math_division"

# Update the active selection / scroll position on the sheet view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F4").Select()
